$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 14.15440837047059
$ws.Range("C2").Value = 4.763913061872104
$ws.Range("D2").Value = 9.106701891442407
$ws.Range("E2").Value = 10.42840865307153
$ws.Range("F2").Value = 45.3110708462298
$ws.Range("I2").Value = 37.03918245434825
$ws.Range("K2").Value = 12.59491124901355
$ws.Range("L2").Value = 10.54869075530701
$ws.Range("M2").Value = 15.46849324139424
# Row 3
$ws.Range("B3").Value = 14.0791922704248
$ws.Range("C3").Value = 4.521650149324073
$ws.Range("D3").Value = 9.096251255057926
$ws.Range("E3").Value = 10.44012375205457
$ws.Range("F3").Value = 44.98296002507977
$ws.Range("I3").Value = 36.8789051741067
$ws.Range("K3").Value = 12.53844627163413
$ws.Range("L3").Value = 10.5569032204261
$ws.Range("M3").Value = 15.4806164093753
# Row 4
$ws.Range("B4").Value = 14.03759974515391
$ws.Range("C4").Value = 4.365059952006747
$ws.Range("D4").Value = 9.089571111249414
$ws.Range("E4").Value = 10.44809934928632
$ws.Range("F4").Value = 44.78664169360824
$ws.Range("I4").Value = 36.78399229748963
$ws.Range("K4").Value = 12.50763355593765
$ws.Range("L4").Value = 10.56341049827718
$ws.Range("M4").Value = 15.49112978020073
# Row 5
$ws.Range("B5").Value = 14.0218222598864
$ws.Range("C5").Value = 4.299305600786987
$ws.Range("D5").Value = 9.086782373268077
$ws.Range("E5").Value = 10.45154654763921
$ws.Range("F5").Value = 44.70798604298087
$ws.Range("I5").Value = 36.74621056018929
$ws.Range("K5").Value = 12.4960588060891
$ws.Range("L5").Value = 10.56643086163664
$ws.Range("M5").Value = 15.49618636401402
# Row 6
$ws.Range("B6").Value = 14.01927366673491
$ws.Range("C6").Value = 4.288270775520095
$ws.Range("D6").Value = 9.086315278974947
$ws.Range("E6").Value = 10.4521308645325
$ws.Range("F6").Value = 44.69500808917498
$ws.Range("I6").Value = 36.7399915077231
$ws.Range("K6").Value = 12.49419643437469
$ws.Range("L6").Value = 10.56695465912577
$ws.Range("M6").Value = 15.49707265752905
# Row 7
$ws.Range("B7").Value = 14.03738219826472
$ws.Range("C7").Value = 4.364180988735125
$ws.Range("D7").Value = 9.089533771106124
$ws.Range("E7").Value = 10.44814504102365
$ws.Range("F7").Value = 44.78557540022165
$ws.Range("I7").Value = 36.78347911100388
$ws.Range("K7").Value = 12.50747346556257
$ws.Range("L7").Value = 10.5634497392368
$ws.Range("M7").Value = 15.49119484772509
# Row 8
$ws.Range("B8").Value = 14.12753221161255
$ws.Range("C8").Value = 4.682023433956563
$ws.Range("D8").Value = 9.10315300927323
$ws.Range("E8").Value = 10.43228580670784
$ws.Range("F8").Value = 45.19690295623766
$ws.Range("I8").Value = 36.98320299637237
$ws.Range("K8").Value = 12.57464913939032
$ws.Range("L8").Value = 10.55121857141863
$ws.Range("M8").Value = 15.47203650666447
# Row 9
$ws.Range("B9").Value = 14.33987600051756
$ws.Range("C9").Value = 5.242122731514793
$ws.Range("D9").Value = 9.127783482380663
$ws.Range("E9").Value = 10.40738069830967
$ws.Range("F9").Value = 46.04191804463395
$ws.Range("I9").Value = 37.40185415276847
$ws.Range("K9").Value = 12.73639499215225
$ws.Range("L9").Value = 10.53884257471802
$ws.Range("M9").Value = 15.45879720280748
# Row 10
$ws.Range("B10").Value = 14.51627967081287
$ws.Range("C10").Value = 5.614117291717548
$ws.Range("D10").Value = 9.144632050779272
$ws.Range("E10").Value = 10.39284100955345
$ws.Range("F10").Value = 46.68281271699456
$ws.Range("I10").Value = 37.72481677942577
$ws.Range("K10").Value = 12.87264201442691
$ws.Range("L10").Value = 10.5368069972074
$ws.Range("M10").Value = 15.46385355267207
# Row 11
$ws.Range("B11").Value = 14.60065147037846
$ws.Range("C11").Value = 5.774634069343588
$ws.Range("D11").Value = 9.15202997766033
$ws.Range("E11").Value = 10.3870385955089
$ws.Range("F11").Value = 46.97798147035282
$ws.Range("I11").Value = 37.87483179524324
$ws.Range("K11").Value = 12.93818980228964
$ws.Range("L11").Value = 10.53740751975712
$ws.Range("M11").Value = 15.46934778368549
# Row 12
$ws.Range("B12").Value = 14.63316619133193
$ws.Range("C12").Value = 5.834158085481592
$ws.Range("D12").Value = 9.154793340403089
$ws.Range("E12").Value = 10.38495776642915
$ws.Range("F12").Value = 47.09021013130694
$ws.Range("I12").Value = 37.93206140280967
$ws.Range("K12").Value = 12.96350435470608
$ws.Range("L12").Value = 10.53785377454323
$ws.Range("M12").Value = 15.4718857222854
# Row 13
$ws.Range("B13").Value = 14.62613889899451
$ws.Range("C13").Value = 5.821394680119434
$ws.Range("D13").Value = 9.154199889920042
$ws.Range("E13").Value = 10.38540073744474
$ws.Range("F13").Value = 47.06602052323456
$ws.Range("I13").Value = 37.91971758506336
$ws.Range("K13").Value = 12.95803081562084
$ws.Range("L13").Value = 10.53774794361824
$ws.Range("M13").Value = 15.47131882272749
# Row 14
$ws.Range("B14").Value = 14.60331531911282
$ws.Range("C14").Value = 5.779556441998404
$ws.Range("D14").Value = 9.152258082706819
$ws.Range("E14").Value = 10.38686507303639
$ws.Range("F14").Value = 46.98720586009082
$ws.Range("I14").Value = 37.87953178447241
$ws.Range("K14").Value = 12.94026268330025
$ws.Range("L14").Value = 10.53743985141926
$ws.Range("M14").Value = 15.46954742659599
# Row 15
$ws.Range("B15").Value = 14.58940793441197
$ws.Range("C15").Value = 5.753764988956625
$ws.Range("D15").Value = 9.1510637193532
$ws.Range("E15").Value = 10.38777717275207
$ws.Range("F15").Value = 46.93898683321471
$ws.Range("I15").Value = 37.854971097348
$ws.Range("K15").Value = 12.9294427828912
$ws.Range("L15").Value = 10.53727961569167
$ws.Range("M15").Value = 15.46852189954231
# Row 16
$ws.Range("B16").Value = 14.51084639430332
$ws.Range("C16").Value = 5.603451196724584
$ws.Range("D16").Value = 9.144143237162243
$ws.Range("E16").Value = 10.39323651921293
$ws.Range("F16").Value = 46.66358999514335
$ws.Range("I16").Value = 37.71507330785375
$ws.Range("K16").Value = 12.86842842138652
$ws.Range("L16").Value = 10.53679841447851
$ws.Range("M16").Value = 15.46355860806073
# Row 17
$ws.Range("B17").Value = 14.46368822684994
$ws.Range("C17").Value = 5.509002152256747
$ws.Range("D17").Value = 9.139829657106862
$ws.Range("E17").Value = 10.39679334360048
$ws.Range("F17").Value = 46.49552391494291
$ws.Range("I17").Value = 37.63002763727689
$ws.Range("K17").Value = 12.83189804694166
$ws.Range("L17").Value = 10.53689374714501
$ws.Range("M17").Value = 15.46133060240832
# Row 18
$ws.Range("B18").Value = 14.43695411838067
$ws.Range("C18").Value = 5.453858497625371
$ws.Range("D18").Value = 9.1373234936208
$ws.Range("E18").Value = 10.39891556088552
$ws.Range("F18").Value = 46.39920175177519
$ws.Range("I18").Value = 37.58140463179215
$ws.Range("K18").Value = 12.81122388059492
$ws.Range("L18").Value = 10.53709230795124
$ws.Range("M18").Value = 15.46034988691204
# Row 19
$ws.Range("B19").Value = 14.42797023529335
$ws.Range("C19").Value = 5.435047423907307
$ws.Range("D19").Value = 9.136470629673358
$ws.Range("E19").Value = 10.39964724282501
$ws.Range("F19").Value = 46.36664998123812
$ws.Range("I19").Value = 37.56499273912066
$ws.Range("K19").Value = 12.80428246278393
$ws.Range("L19").Value = 10.53718424139291
$ws.Range("M19").Value = 15.46006954513192
# Row 20
$ws.Range("B20").Value = 14.46866812916586
$ws.Range("C20").Value = 5.519141252245694
$ws.Range("D20").Value = 9.140291441583443
$ws.Range("E20").Value = 10.39640680601319
$ws.Range("F20").Value = 46.51337960733571
$ws.Range("I20").Value = 37.63905073439256
$ws.Range("K20").Value = 12.83575201754995
$ws.Range("L20").Value = 10.53686872805991
$ws.Range("M20").Value = 15.46153666136146
# Row 21
$ws.Range("B21").Value = 14.61000405549243
$ws.Range("C21").Value = 5.791879603431341
$ws.Range("D21").Value = 9.152829469746044
$ws.Range("E21").Value = 10.38643180510189
$ws.Range("F21").Value = 47.01034382974397
$ws.Range("I21").Value = 37.89132405971898
$ws.Range("K21").Value = 12.94546839880298
$ws.Range("L21").Value = 10.53752441173486
$ws.Range("M21").Value = 15.47005533177784
# Row 22
$ws.Range("B22").Value = 14.70565438827419
$ws.Range("C22").Value = 5.96278213881088
$ws.Range("D22").Value = 9.160802003894872
$ws.Range("E22").Value = 10.38059099726842
$ws.Range("F22").Value = 47.33775756086185
$ws.Range("I22").Value = 38.05864871053942
$ws.Range("K22").Value = 13.02003815903166
$ws.Range("L22").Value = 10.53922819678813
$ws.Range("M22").Value = 15.4782877131476
# Row 23
$ws.Range("B23").Value = 14.65431356334513
$ws.Range("C23").Value = 5.872242737330566
$ws.Range("D23").Value = 9.156567121122469
$ws.Range("E23").Value = 10.38364637400682
$ws.Range("F23").Value = 47.16279314071596
$ws.Range("I23").Value = 37.96912815309087
$ws.Range("K23").Value = 12.97998372391875
$ws.Range("L23").Value = 10.53820240897092
$ws.Range("M23").Value = 15.47365080236748
# Row 24
$ws.Range("B24").Value = 14.4664155365584
$ws.Range("C24").Value = 5.514559992282314
$ws.Range("D24").Value = 9.140082750533201
$ws.Range("E24").Value = 10.39658131870152
$ws.Range("F24").Value = 46.50530610530262
$ws.Range("I24").Value = 37.63497054660503
$ws.Range("K24").Value = 12.8340086153907
$ws.Range("L24").Value = 10.53687959137272
$ws.Range("M24").Value = 15.46144256690057
# Row 25
$ws.Range("B25").Value = 14.27875025774953
$ws.Range("C25").Value = 5.097483332068819
$ws.Range("D25").Value = 9.121340682534164
$ws.Range("E25").Value = 10.41345687919228
$ws.Range("F25").Value = 45.80954633866862
$ws.Range("I25").Value = 37.2858152377402
$ws.Range("K25").Value = 12.68951132431454
$ws.Range("L25").Value = 10.5409495322952
$ws.Range("M25").Value = 15.45977802025113
